# Adds the next day's COVID-19 patient-condition record to the
# "Condicion_Pacientes" table on Hoja1 (row 87: 07/06/2020 data) and
# moves the active selection to the newly added row, mirroring the
# manual "Actualizar" workflow used for this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Last existing data row of the table (row 86, A86:F86) supplies the
# formatting (date format on column A, centered numbers on B:F) for
# the new row.
$lastRow = $lo.ListRows.Item($lo.ListRows.Count).Range
$lastRow.Copy()

# Grow the table by one row; the new ListRow's range is A87:F87 and the
# table/autoFilter ref expands from A1:F86 to A1:F87 automatically.
$newListRow = $lo.ListRows.Add()
$newRange = $newListRow.Range
$newRange.PasteSpecial(-4122)

# New day's figures.
$newRange.Cells.Item(1, 1).Value2 = 43989
$newRange.Cells.Item(1, 2).Value2 = 647
$newRange.Cells.Item(1, 3).Value2 = 172
$newRange.Cells.Item(1, 4).Value2 = 430
$newRange.Cells.Item(1, 5).Value2 = 143
$newRange.Cells.Item(1, 6).Value2 = 47

# Match the author's saved view state: selection on the new last cell.
$ws.Activate()
$newRange.Cells.Item(1, 6).Select()
